$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for cetsid "ENF.CONT.COEN.ENFE.PR" (Enforcing contracts: Enforcement
# fees (% of claim)) was removed from the dataset, so delete its entire row.
# This is row 102 (A102 = "ENF.CONT.COEN.ENFE.PR").
[void]$ws.Rows.Item(102).Delete()

# The worksheet carries a stale hidden _FilterDatabase defined name that still
# spans the old used range; keep it in sync with the new (smaller) extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$M`$135"
    }
}

# Refresh the active selection in the lower (scrolled/frozen) pane.
[void]$ws.Range("C63").Select()
